$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '20.498.26'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +2.35%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.475.73'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +4.11%  '
$ws.Range('E4').Value = '  +0.84%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.9636'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -3.58%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '276.64'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.16%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3646'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -1.11%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3044'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -2.14%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '39.75'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.29%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.049'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.84%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.06596'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.24%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.002'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.42%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '18.13'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +2.84%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.453'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.77%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.162'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.39%  '
$ws.Range('E16').Value = '  +0.91%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.478.40'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +4.64%  '
$ws.Range('E18').Value = '  +3.74%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.9695'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -3.07%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '68.88'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -3.19%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.461'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -2.55%  '
$ws.Range('E22').Value = '  -2.26%  '
$ws.Range('E23').Value = '  -0.59%  '
$ws.Range('E24').Value = '  +0.35%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '20.526.42'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +2.49%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '141.48'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +6.16%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.127'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -6.77%  '
$ws.Range('E28').Value = '  -0.55%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.632.79'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +3.89%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '113.32'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +2.77%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.902'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.27%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.951'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -5.07%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.8100'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.57%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.07866'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.13%  '
$ws.Range('B35').Value = 'WEMIXTOKEN'
$ws.Range('C35').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.518'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +3.92%  '
$ws.Range('B36').Value = 'TrustWalletToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.248'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +13.09%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.05724'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.80%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.725'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -3.55%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '7.741'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -4.01%  '
$ws.Range('E40').Value = '  -0.67%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.9629'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -3.66%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '10.39'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.81%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.1873'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.16%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.5266'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.75%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.502'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.90%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '12.04'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -2.27%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '116.65'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.67%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.5158'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.56%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.771'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.13%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06442'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +3.89%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.9896'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.07%  '
